$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G4").Value = 'Dr. Hend Farid, Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Amal Awwad'
$ws.Range("G5").Value = 'Dr. Amal Awwad, D Wessam Atef, Dr. Sara Nabil, Dr. Nourhan Mohammad'
$ws.Range("G6").Value = 'Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad'
$ws.Range("G8").Value = 'Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G9").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G10").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G15").Value = 'Dr. Ahmad Mostafa, Dr. Marian Samir, Dr. Afaf Abdallah, Dr. Nourham Mostafa, Dr. Nourhan Mohammad'
$ws.Range("G16").Value = 'Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Rada Rabea, Dr. Walaa Ghanima'
$ws.Range("G17").Value = 'Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Nardine, Dr. Monica, Dr. Youstina Magdy'
$ws.Range("G18").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Remon, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida'
$ws.Range("G19").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G21").Value = 'Dr. Hend Farid, Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Amal Awwad'
$ws.Range("G22").Value = 'Dr. Amal Awwad, D Wessam Atef, Dr. Sara Nabil, Dr. Nourhan Mohammad'
$ws.Range("G23").Value = 'Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad'
$ws.Range("G25").Value = 'Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G26").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G27").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G30").Value = 'Dr. Sarah Mahdy, Dr. Mariam Gamal Sanad'
$ws.Range("G32").Value = 'Dr. Ahmad Mostafa, Dr. Marian Samir, Dr. Afaf Abdallah, Dr. Nourham Mostafa, Dr. Nourhan Mohammad'
$ws.Range("G33").Value = 'Dr. Marian Samir, Dr. Manarst Al-Eslam, Dr. Rada Rabea, Dr. Walaa Ghanima'
$ws.Range("G34").Value = 'Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Nardine, Dr. Monica, Dr. Youstina Magdy'
$ws.Range("G35").Value = 'Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Remon'
$ws.Range("G36").Value = 'Dr. Eman Tantawi, Administrator, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G37").Value = 'Administrator, Dr. Kerelos Zareef, Dr. Nada Mohammad'
$ws.Range("G40").Value = 'Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad'
$ws.Range("G43").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G44").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G45").Value = 'Dr. Mohammad Safwat, Administrator, Dr. Rania Ahmad Youssef'
$ws.Range("G48").Value = 'Dr. Afaf Abdallah, Dr. Aya Alaa-Eldein, Dr. Marian Samir'
$ws.Range("G49").Value = 'Dr. Nourham Mostafa, Dr. Aya Alaa-Eldein'
$ws.Range("G50").Value = 'Dr. Manarst Al-Eslam, Dr. Nancy Abd Al-Shafy, Dr. Aya Alaa-Eldein'
$ws.Range("G51").Value = 'Dr. Naema Gomaa, Dr. Yasmin, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range("G52").Value = 'Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Yasmin'
$ws.Range("G53").Value = 'Dr. Eman Tantawi, Administrator, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G54").Value = 'Administrator, Dr. Kerelos Zareef, Dr. Nada Mohammad'
$ws.Range("G57").Value = 'Dr. Mai Mustafa, Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Maryam Ahmad'
$ws.Range("G60").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G61").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G62").Value = 'Dr. Mohammad Safwat, Administrator, Dr. Rania Ahmad Youssef'
$ws.Range("G65").Value = 'Dr. Afaf Abdallah, Dr. Aya Alaa-Eldein, Dr. Marian Samir'
$ws.Range("G66").Value = 'Dr. Nourham Mostafa, Dr. Aya Alaa-Eldein'
$ws.Range("G67").Value = 'Dr. Manarst Al-Eslam, Dr. Nancy Abd Al-Shafy, Dr. Aya Alaa-Eldein'
$ws.Range("G68").Value = 'Dr. Naema Gomaa, Dr. Yasmin, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Monica, Dr. Wafaa Ebida'
$ws.Range("G69").Value = 'Dr. Naema Gomaa, Dr. Salma Hassan, Dr. Abdullah El-Agrody, Dr. Remon, Dr. Yasmin'
$ws.Range("G70").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G72").Value = 'Dr. Mariam Nour El-Din, Dr. Safa Hany, Dr. Shimaa Ashraf, D Wessam Atef, Dr. Omnia Mohammad'
$ws.Range("G73").Value = 'Dr. Hend Farid, Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Amal Awwad'
$ws.Range("G74").Value = 'Dr. Sara Nabil, Dr. Aya Saeed, D Wessam Atef, Dr. Omnia Mohammad, Dr. Amal Awwad'
$ws.Range("G75").Value = 'Dr. Eman M. Elsaid, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad'
$ws.Range("G76").Value = 'Dr. Mohammad Safwat, Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range("G77").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range("G78").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Al-Shimaa Khaled'
$ws.Range("G79").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G83").Value = 'Dr. Afaf Abdallah, Dr. Youstina Ibrahim, Dr. Marian Samir'
$ws.Range("G84").Value = 'Dr. Manarst Al-Eslam, Dr. Nancy Abd Al-Shafy, Dr. Aya Alaa-Eldein, Dr. Marian Samir'
$ws.Range("G85").Value = 'Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nahla, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida'
$ws.Range("G86").Value = 'Dr. Marina Atef, Dr. Yasmin, Dr. Salma Hassan, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Remon, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida'
$ws.Range("G87").Value = 'Dr. Eman Tantawi, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G88").Value = 'Dr. Fatma Elhady, Dr. Nada Mohammad'
$ws.Range("G89").Value = 'Dr. Hend Farid, Dr. Mariam Nour El-Din, Dr. Aya Saeed, Dr. Shimaa Ashraf, Dr. Amal Awwad'
$ws.Range("G90").Value = 'Dr. Sara Nabil, Dr. Aya Saeed, D Wessam Atef, Dr. Omnia Mohammad, Dr. Amal Awwad'
$ws.Range("G91").Value = 'Dr. Eman M. Elsaid, Dr. Merna Said, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya, Dr. Maryam Ahmad'
$ws.Range("G92").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range("G93").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G94").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G98").Value = 'Dr. Nourhan Hosni, Dr. Afaf Abdallah, Dr. Nourhan Mohammad, Dr. Walaa Ghanima'
$ws.Range("G100").Value = 'Dr. Afaf Abdallah, Dr. Nourhan Mohammad, Dr. Aya Alaa-Eldein'
$ws.Range("G101").Value = 'Dr. Manarst Al-Eslam, Dr. Nancy Abd Al-Shafy, Dr. Aya Alaa-Eldein, Dr. Marian Samir'
$ws.Range("G102").Value = 'Dr. Yasmin, Dr. Neveen Nashaat, Dr. Nahla, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Maryam Ashraf, Dr. Wafaa Ebida'
$ws.Range("G103").Value = 'Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Gehad Salah, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Maryam Ashraf'
$ws.Range("G104").Value = 'Dr. Eman Tantawi, Administrator, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G106").Value = 'Dr. Amal Awwad, D Wessam Atef, Dr. Sara Nabil, Dr. Nourhan Mohammad'
$ws.Range("G107").Value = 'Dr. Amal Awwad, Dr. Nourhan Mohammad, Dr. Sara Nabil'
$ws.Range("G108").Value = 'Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya'
$ws.Range("G111").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G112").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G115").Value = 'Dr. Sarah Mahdy, Dr. Mariam Gamal Sanad'
$ws.Range("G116").Value = 'Dr. Afaf Abdallah, Dr. Nourham Mostafa, Dr. Enas Omran'
$ws.Range("G117").Value = 'Dr. Taqwa Mohammad, Dr. Amr Saeed, Dr. Enas Omran'
$ws.Range("G119").Value = 'Dr. Marina Atef, Dr. Remon, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nahla, Dr. Marina Sorial, Dr. Nardine, Dr. Aya Hanafy, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range("G120").Value = 'Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Gehad Salah, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Maryam Ashraf'
$ws.Range("G121").Value = 'Dr. Eman Tantawi, Administrator, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G123").Value = 'Dr. Amal Awwad, D Wessam Atef, Dr. Sara Nabil, Dr. Nourhan Mohammad'
$ws.Range("G124").Value = 'Dr. Amal Awwad, Dr. Nourhan Mohammad, Dr. Sara Nabil'
$ws.Range("G128").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Al-Shimaa Khaled'
$ws.Range("G129").Value = 'Administrator, Dr. Al-Shimaa Khaled, Dr. Rania Ahmad Youssef'
$ws.Range("G132").Value = 'Dr. Sarah Mahdy, Dr. Mariam Gamal Sanad'
$ws.Range("G133").Value = 'Dr. Afaf Abdallah, Dr. Nourham Mostafa, Dr. Enas Omran'
$ws.Range("G134").Value = 'Dr. Taqwa Mohammad, Dr. Amr Saeed, Dr. Enas Omran'
$ws.Range("G136").Value = 'Dr. Marina Atef, Dr. Remon, Dr. Shorok Mohammad, Dr. Neveen Nashaat, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Nahla, Dr. Marina Sorial, Dr. Nardine, Dr. Aya Hanafy, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Range("G137").Value = 'Dr. Yassmen Ahmad, Dr. Salma Hassan, Dr. Gehad Salah, Dr. Neveen Nashaat, Dr. Abdullah El-Agrody, Dr. Aya Emad, Dr. Ola Abd Al-Fattah, Dr. Aya Hanafy, Dr. Maryam Ashraf'
$ws.Range("G138").Value = 'Dr. Eman Tantawi, Administrator, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Rana Abo-Zaid'
$ws.Range("G142").Value = 'Dr. Amany Raafat, Dr. Madeha Saeed, Dr. Merna Said, Dr. Basma Hamed, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Eman M. Abo-Sakaya'
$ws.Range("G144").Value = 'Dr. Mayar Ahmad Embaby, Nourhan Mamdouh Hassan, Dr. Mohammad Safwat, Dr. Mariam Toma Gerges'
$ws.Range("G145").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Al-Shimaa Khaled'
$ws.Range("G146").Value = 'Dr. Mayar Ahmad Embaby, Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef'
$ws.Range("G148").Value = 'Dr. Sarah Mahdy, Dr. Mariam Gamal Sanad'
$ws.Range("G150").Value = 'Dr. Afaf Abdallah, Dr. Youstina Ibrahim, Dr. Marian Samir'
$ws.Range("G151").Value = 'Administrator, Dr. Rada Rabea, Dr. Marian Samir, Dr. Hana Amr, Dr. Nourhan Mohammad'
